{"js": "// \"Getting rid of AR\" \u2014 all occurrences of the \"AR\" abbreviation (as a\n// stand-alone word, including inside the ${AR} placeholder) are replaced\n// with \"Recommendation\" / \"recommendation\", per the commit message:\n//   \"All ARs are replaced with recommendations\"\n//\n// The document contains four \"AR\" substrings in total:\n//   1) the title run \"AR\" (-> \"Recommendation\")\n//   2) the title placeholder \"${AR}\" (-> \"${REC}\")\n//   3) \"ARC Number\" inside a table (unrelated - must NOT change)\n//   4) \"...for this AR is...\" (-> \"...for this recommendation is...\")\n//\n// We target (1)+(2) together (\"AR ${AR}:\" is unique in the document) and\n// (4) separately (\"this AR is\" is unique in the document) with literal,\n// case-sensitive searches so the unrelated \"ARC Number\" table cell is left\n// untouched.\n\nconst body = context.document.body;\n\n// --- Title: \"AR ${AR}:\" -> \"Recommendation ${REC}:\" ------------------\nconst titleResults = body.search(\"AR ${AR}:\", { matchCase: true, matchWholeWord: false });\ntitleResults.load(\"text\");\nawait context.sync();\n\n// --- Body sentence: \"this AR is\" -> \"this recommendation is\" ----------\nconst sentenceResults = body.search(\"this AR is\", { matchCase: true, matchWholeWord: false });\nsentenceResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < titleResults.items.length; i++) {\n  titleResults.items[i].insertText(\"Recommendation ${REC}:\", Word.InsertLocation.replace);\n}\n\nfor (let i = 0; i < sentenceResults.items.length; i++) {\n  sentenceResults.items[i].insertText(\"this recommendation is\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# \"Getting rid of AR\" \u2014 all occurrences of the \"AR\" abbreviation (as a\n# stand-alone word, including inside the ${AR} placeholder) are replaced\n# with \"Recommendation\" / \"recommendation\", per the commit message:\n#   \"All ARs are replaced with recommendations\"\n#\n# The document contains four \"AR\" substrings in total:\n#   1) the title run \"AR\"              -> \"Recommendation\"\n#   2) the title placeholder \"${AR}\"   -> \"${REC}\"\n#   3) \"ARC Number\" inside a table     -> unrelated, must NOT change\n#   4) \"...for this AR is...\"          -> \"...for this recommendation is...\"\n#\n# We target (1)+(2) together (\"AR ${AR}:\" is unique in the document) and\n# (4) separately (\"this AR is\" is unique in the document) with literal,\n# case-sensitive Find/Replace so the unrelated \"ARC Number\" table cell is\n# left untouched.\n#\n# NOTE: PowerShell double-quoted strings interpolate \"${AR}\" as a variable\n# reference, so single-quoted strings are used for every literal that\n# contains a \"${...}\" placeholder.\n\n$d = $word.ActiveDocument\n\n# wdFindWrap: wdFindContinue = 1 ; wdReplace: wdReplaceOne = 1, wdReplaceAll = 2\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n# --- Title: \"AR ${AR}:\" -> \"Recommendation ${REC}:\" --------------------\n$range1 = $d.Content\n$found1 = $range1.Find.Execute('AR ${AR}:', $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, 'Recommendation ${REC}:', $wdReplaceAll)\n\n# --- Body sentence: \"this AR is\" -> \"this recommendation is\" -----------\n$range2 = $d.Content\n$found2 = $range2.Find.Execute('this AR is', $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, 'this recommendation is', $wdReplaceAll)\n\nWrite-Output \"title replaced: $found1; sentence replaced: $found2\"\n"}
